# Add data for 2024-03-10
# Applies the updated 2024 year-to-date crime counts (column K = year 2024)
# across the "Citywide Totals", "By Neighborhood" and individual
# neighborhood worksheets, reflecting one additional day of reported
# incidents (2024-03-10).

$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet='Citywide Totals'; Cell='K2'; Value=1322},
    @{Sheet='Citywide Totals'; Cell='K3'; Value=1253},
    @{Sheet='Citywide Totals'; Cell='I4'; Value=1787},
    @{Sheet='Citywide Totals'; Cell='K4'; Value=272},
    @{Sheet='Citywide Totals'; Cell='K5'; Value=79},
    @{Sheet='Citywide Totals'; Cell='K6'; Value=1600},
    @{Sheet='Citywide Totals'; Cell='I7'; Value=26240},
    @{Sheet='Citywide Totals'; Cell='K7'; Value=4526},
    @{Sheet='Austin'; Cell='K3'; Value=78},
    @{Sheet='Austin'; Cell='K7'; Value=263},
    @{Sheet='South Chicago'; Cell='K4'; Value=4},
    @{Sheet='South Chicago'; Cell='K7'; Value=90},
    @{Sheet='Garfield Park'; Cell='K3'; Value=72},
    @{Sheet='Garfield Park'; Cell='K7'; Value=186},
    @{Sheet='Grand Crossing'; Cell='K3'; Value=52},
    @{Sheet='Grand Crossing'; Cell='K6'; Value=56},
    @{Sheet='Grand Crossing'; Cell='K7'; Value=152},
    @{Sheet='Woodlawn'; Cell='K3'; Value=30},
    @{Sheet='Woodlawn'; Cell='K6'; Value=25},
    @{Sheet='By Neighborhood'; Cell='K6'; Value=35},
    @{Sheet='By Neighborhood'; Cell='K7'; Value=131},
    @{Sheet='By Neighborhood'; Cell='K8'; Value=263},
    @{Sheet='By Neighborhood'; Cell='K11'; Value=96},
    @{Sheet='By Neighborhood'; Cell='K14'; Value=28},
    @{Sheet='By Neighborhood'; Cell='K18'; Value=35},
    @{Sheet='By Neighborhood'; Cell='K19'; Value=121},
    @{Sheet='By Neighborhood'; Cell='K20'; Value=110},
    @{Sheet='By Neighborhood'; Cell='K22'; Value=11},
    @{Sheet='By Neighborhood'; Cell='K23'; Value=42},
    @{Sheet='By Neighborhood'; Cell='K27'; Value=54},
    @{Sheet='By Neighborhood'; Cell='K29'; Value=212},
    @{Sheet='By Neighborhood'; Cell='K32'; Value=9},
    @{Sheet='By Neighborhood'; Cell='K33'; Value=186},
    @{Sheet='By Neighborhood'; Cell='K36'; Value=50},
    @{Sheet='By Neighborhood'; Cell='K37'; Value=152},
    @{Sheet='By Neighborhood'; Cell='K41'; Value=46},
    @{Sheet='By Neighborhood'; Cell='K42'; Value=153},
    @{Sheet='By Neighborhood'; Cell='K43'; Value=44},
    @{Sheet='By Neighborhood'; Cell='J45'; Value=41},
    @{Sheet='By Neighborhood'; Cell='K47'; Value=34},
    @{Sheet='By Neighborhood'; Cell='K48'; Value=49},
    @{Sheet='By Neighborhood'; Cell='K52'; Value=123},
    @{Sheet='By Neighborhood'; Cell='K54'; Value=79},
    @{Sheet='By Neighborhood'; Cell='K55'; Value=47},
    @{Sheet='By Neighborhood'; Cell='K56'; Value=7},
    @{Sheet='By Neighborhood'; Cell='I63'; Value=196},
    @{Sheet='By Neighborhood'; Cell='J63'; Value=89},
    @{Sheet='By Neighborhood'; Cell='K63'; Value=17},
    @{Sheet='By Neighborhood'; Cell='K64'; Value=30},
    @{Sheet='By Neighborhood'; Cell='K66'; Value=19},
    @{Sheet='By Neighborhood'; Cell='K67'; Value=183},
    @{Sheet='By Neighborhood'; Cell='K76'; Value=62},
    @{Sheet='By Neighborhood'; Cell='K78'; Value=64},
    @{Sheet='By Neighborhood'; Cell='K79'; Value=121},
    @{Sheet='By Neighborhood'; Cell='K83'; Value=90},
    @{Sheet='By Neighborhood'; Cell='K85'; Value=229},
    @{Sheet='By Neighborhood'; Cell='K88'; Value=57},
    @{Sheet='By Neighborhood'; Cell='K91'; Value=47},
    @{Sheet='By Neighborhood'; Cell='K97'; Value=37},
    @{Sheet='By Neighborhood'; Cell='K98'; Value=28},
    @{Sheet='By Neighborhood'; Cell='I101'; Value=26240},
    @{Sheet='By Neighborhood'; Cell='K101'; Value=4526},
    @{Sheet='North Lawndale'; Cell='K5'; Value=4},
    @{Sheet='North Lawndale'; Cell='K7'; Value=183},
    @{Sheet='Loop'; Cell='K2'; Value=17},
    @{Sheet='Loop'; Cell='K6'; Value=31},
    @{Sheet='Loop'; Cell='K7'; Value=79},
    @{Sheet='Englewood'; Cell='K3'; Value=67},
    @{Sheet='Englewood'; Cell='K4'; Value=9},
    @{Sheet='Englewood'; Cell='K6'; Value=76},
    @{Sheet='Englewood'; Cell='K7'; Value=212},
    @{Sheet='Lake View'; Cell='K3'; Value=9},
    @{Sheet='Lake View'; Cell='K7'; Value=49},
    @{Sheet='Chatham'; Cell='K2'; Value=39},
    @{Sheet='Chatham'; Cell='K6'; Value=36},
    @{Sheet='Chatham'; Cell='K7'; Value=121},
    @{Sheet='River North'; Cell='K3'; Value=12},
    @{Sheet='River North'; Cell='K6'; Value=32},
    @{Sheet='River North'; Cell='K7'; Value=62},
    @{Sheet='Bridgeport'; Cell='K2'; Value=12},
    @{Sheet='Bridgeport'; Cell='K6'; Value=10},
    @{Sheet='Bridgeport'; Cell='K7'; Value=28},
    @{Sheet='Ashburn'; Cell='K2'; Value=12},
    @{Sheet='Ashburn'; Cell='K6'; Value=11},
    @{Sheet='Ashburn'; Cell='K7'; Value=35},
    @{Sheet='Hermosa'; Cell='K2'; Value=12},
    @{Sheet='Hermosa'; Cell='K5'; Value=1},
    @{Sheet='Hermosa'; Cell='K7'; Value=46},
    @{Sheet='Humboldt Park'; Cell='K2'; Value=39},
    @{Sheet='Humboldt Park'; Cell='K7'; Value=153},
    @{Sheet='Rogers Park'; Cell='K4'; Value=4},
    @{Sheet='Rogers Park'; Cell='K7'; Value=64},
    @{Sheet='Lower West Side'; Cell='K3'; Value=9},
    @{Sheet='Lower West Side'; Cell='K7'; Value=47},
    @{Sheet='Douglas'; Cell='K3'; Value=11},
    @{Sheet='Douglas'; Cell='K4'; Value=6},
    @{Sheet='Douglas'; Cell='K7'; Value=42},
    @{Sheet='Washington Park'; Cell='K3'; Value=17},
    @{Sheet='Washington Park'; Cell='K7'; Value=47},
    @{Sheet='Roseland'; Cell='K2'; Value=45},
    @{Sheet='Roseland'; Cell='K3'; Value=39},
    @{Sheet='Roseland'; Cell='K7'; Value=121},
    @{Sheet='Near South Side'; Cell='K3'; Value=11},
    @{Sheet='Near South Side'; Cell='K7'; Value=30},
    @{Sheet='Chicago Lawn'; Cell='K2'; Value=32},
    @{Sheet='Chicago Lawn'; Cell='K6'; Value=42},
    @{Sheet='Chicago Lawn'; Cell='K7'; Value=110},
    @{Sheet='Calumet Heights'; Cell='K5'; Value=1},
    @{Sheet='Calumet Heights'; Cell='K7'; Value=35},
    @{Sheet='Grand Boulevard'; Cell='K3'; Value=19},
    @{Sheet='Grand Boulevard'; Cell='K7'; Value=50},
    @{Sheet='Auburn Gresham'; Cell='K3'; Value=39},
    @{Sheet='Auburn Gresham'; Cell='K7'; Value=131},
    @{Sheet='Kenwood'; Cell='K6'; Value=10},
    @{Sheet='Kenwood'; Cell='K7'; Value=34},
    @{Sheet='Wicker Park'; Cell='K6'; Value=21},
    @{Sheet='Wicker Park'; Cell='K7'; Value=28},
    @{Sheet='North Center'; Cell='K6'; Value=9},
    @{Sheet='North Center'; Cell='K7'; Value=19},
    @{Sheet='Belmont Cragin'; Cell='K3'; Value=23},
    @{Sheet='Belmont Cragin'; Cell='K7'; Value=96},
    @{Sheet='West Town'; Cell='K6'; Value=24},
    @{Sheet='West Town'; Cell='K7'; Value=37},
    @{Sheet='United Center'; Cell='K6'; Value=31},
    @{Sheet='United Center'; Cell='K7'; Value=57},
    @{Sheet='Galewood'; Cell='K3'; Value=3},
    @{Sheet='Galewood'; Cell='K7'; Value=9},
    @{Sheet='Edgewater'; Cell='K2'; Value=19},
    @{Sheet='Edgewater'; Cell='K7'; Value=54},
    @{Sheet='Hyde Park'; Cell='K2'; Value=7},
    @{Sheet='Hyde Park'; Cell='K3'; Value=15},
    @{Sheet='Hyde Park'; Cell='K7'; Value=44},
    @{Sheet='South Shore'; Cell='K2'; Value=85},
    @{Sheet='South Shore'; Cell='K6'; Value=53},
    @{Sheet='South Shore'; Cell='K7'; Value=229},
    @{Sheet='Clearing'; Cell='K2'; Value=4},
    @{Sheet='Clearing'; Cell='K7'; Value=11},
    @{Sheet='Jackson Park'; Cell='J3'; Value=11},
    @{Sheet='Jackson Park'; Cell='J7'; Value=41},
    @{Sheet='Magnificent Mile'; Cell='K3'; Value=2},
    @{Sheet='Magnificent Mile'; Cell='K7'; Value=7},
    @{Sheet='Little Village'; Cell='K3'; Value=26},
    @{Sheet='Little Village'; Cell='K6'; Value=58},
    @{Sheet='Little Village'; Cell='K7'; Value=123}
)

$sheetCache = @{}

foreach ($change in $changes) {
    $sheetName = $change.Sheet
    if (-not $sheetCache.ContainsKey($sheetName)) {
        $sheetCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $sheetCache[$sheetName]
    $ws.Range($change.Cell).Value = $change.Value
}

Write-Host "Applied" $changes.Count "cell updates across" $sheetCache.Count "worksheets."
